$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "Immutable"
$ws.Range("D1").Value = "अडिग"

$ws.Range("D2").Value = "परिवर्तनशील"
$ws.Range("B2").Value = "Mutable"

$ws.Range("D3").Value = "यह दर्शाता है"
$ws.Range("B3").Value = "demonstrates"

$ws.Range("D4").Value = "स्पष्ट रूप से"
$ws.Range("B4").Value = "explicitly"

$ws.Range("B8").Select() | Out-Null
